$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update recorded sensor values for rows 2-17 (new event timing/measurement data) ---
$ws.Range("B2").Value2 = 2.0005650000000514
$ws.Range("C2").Value2 = 1998.0601435999999
$ws.Range("D2").Value2 = 2000.0619563

$ws.Range("B3").Value2 = 4.0000984000000699
$ws.Range("C3").Value2 = 1998.0604370999999
$ws.Range("D3").Value2 = 2002.0614897

$ws.Range("B4").Value2 = 5.9994818000000123
$ws.Range("C4").Value2 = 1998.0604378999999
$ws.Range("D4").Value2 = 2004.0608731

$ws.Range("B5").Value2 = 7.9998361000000386
$ws.Range("C5").Value2 = 1998.0604373000001
$ws.Range("D5").Value2 = 2006.0612274

$ws.Range("B6").Value2 = 9.9995304000001397
$ws.Range("C6").Value2 = 1998.0604390000001
$ws.Range("D6").Value2 = 2008.0609217000001

$ws.Range("B7").Value2 = 11.99948940000013
$ws.Range("C7").Value2 = 1998.0604389
$ws.Range("D7").Value2 = 2010.0608807000001

$ws.Range("B8").Value2 = 13.999432200000001
$ws.Range("C8").Value2 = 1998.0604394
$ws.Range("D8").Value2 = 2012.0608235

$ws.Range("B9").Value2 = 15.999476899999991
$ws.Range("C9").Value2 = 1998.0604346
$ws.Range("D9").Value2 = 2014.0608682

$ws.Range("B10").Value2 = 17.999506399999973
$ws.Range("C10").Value2 = 1998.0604386
$ws.Range("D10").Value2 = 2016.0608976999999

$ws.Range("B11").Value2 = 19.999389999999948
$ws.Range("C11").Value2 = 1998.0604395
$ws.Range("D11").Value2 = 2018.0607812999999

$ws.Range("B12").Value2 = 21.999400599999944
$ws.Range("C12").Value2 = 1998.0604394
$ws.Range("D12").Value2 = 2020.0607918999999

$ws.Range("B13").Value2 = 23.999452100000099
$ws.Range("C13").Value2 = 1998.0604393000001
$ws.Range("D13").Value2 = 2022.0608434000001

$ws.Range("B14").Value2 = 25.999831900000117
$ws.Range("C14").Value2 = 1998.0604369
$ws.Range("D14").Value2 = 2024.0612232000001

$ws.Range("B15").Value2 = 27.999381900000117
$ws.Range("C15").Value2 = 1998.0604393000001
$ws.Range("D15").Value2 = 2026.0607732000001

$ws.Range("B16").Value2 = 30.000027199999977
$ws.Range("C16").Value2 = 1998.0604392
$ws.Range("D16").Value2 = 2028.0614184999999

$ws.Range("B17").Value2 = 31.999517400000059
$ws.Range("C17").Value2 = 1998.0604390000001
$ws.Range("D17").Value2 = 2030.0609087

# --- Rows 18-22 no longer hold trial data: zero out A-D and clear the E (stimuli label) cell ---
foreach ($r in 18..22) {
    $ws.Range("A$r").Value2 = 0
    $ws.Range("B$r").Value2 = 0
    $ws.Range("C$r").Value2 = 0
    $ws.Range("D$r").Value2 = 0
    $ws.Range("E$r").ClearContents()
}
